$wb = $excel.ActiveWorkbook

# Remember which sheet is currently active so we can restore it afterwards.
$origSheet = $wb.ActiveSheet

# Update the ProjectGroups sheet: rename project group labels
# "Road_complete" -> "Road" and "Subway_complete" -> "Rail"
$ws = $wb.Worksheets.Item("ProjectGroups")
$ws.Range("B2").Value = "Road"
$ws.Range("B3").Value = "Rail"

# Select B2 as the active cell on this sheet (records the saved selection
# in the sheet view), then restore the originally active sheet/tab.
$ws.Activate()
$ws.Range("B2").Select()
$origSheet.Activate()
